# "Update countries & provincias Spain"
#
# The underlying data feed re-sorted a couple of countries (by total case
# count) and refreshed the daily COVID-19 counters plus the "last updated"
# timestamp. In the sheet, each data row keeps the same row number but the
# row that used to hold Israel/Egipto/Ucrania (rows 32-34) now holds
# Ucrania/Israel/Egipto (Ucrania overtook the other two), and similarly
# Montserrat/Islas Malvinas (rows 213-214) swap places. We reproduce this
# by writing the resulting country name + stats directly into each cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 09:22"

# Estados Unidos (row 4) - updated counters
$ws.Range("B4").Value = 5701162
$ws.Range("C4").Value = 231
$ws.Range("D4").Value = 3062755
$ws.Range("E4").Value = 2462065
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 176342

# India (row 6) - updated counters
$ws.Range("B6").Value = 2841337
$ws.Range("C6").Value = 5515
$ws.Range("D6").Value = 2097761
$ws.Range("E6").Value = 689559
$ws.Range("G6").Value = 23
$ws.Range("H6").Value = 54017

# Ucrania moves up, now ranked ahead of Israel and Egipto (rows 32-34)
$ws.Range("A32").Value = "Ucrania"
$ws.Range("B32").Value = 98537
$ws.Range("C32").Value = 2134
$ws.Range("D32").Value = 50441
$ws.Range("E32").Value = 45912
$ws.Range("G32").Value = 40
$ws.Range("H32").Value = 2184

$ws.Range("A33").Value = "Israel"
$ws.Range("B33").Value = 98443
$ws.Range("C33").Value = 474
$ws.Range("D33").Value = 73841
$ws.Range("E33").Value = 23821
$ws.Range("H33").Value = 781

$ws.Range("A34").Value = "Egipto"
$ws.Range("B34").Value = 96914
$ws.Range("D34").Value = 62553
$ws.Range("E34").Value = 29164
$ws.Range("H34").Value = 5197

# Armenia (row 57) - updated counters
$ws.Range("B57").Value = 42319
$ws.Range("C57").Value = 263
$ws.Range("D57").Value = 35476
$ws.Range("E57").Value = 6007
$ws.Range("G57").Value = 3
$ws.Range("H57").Value = 836

# Suiza (row 59) - active/recovered reclassified
$ws.Range("D59").Value = 33800
$ws.Range("E59").Value = 2964

# Australia (row 71) - active/recovered reclassified
$ws.Range("D71").Value = 17854
$ws.Range("E71").Value = 5919

# Hungria (row 108) - updated counters
$ws.Range("B108").Value = 5046
$ws.Range("C108").Value = 44
$ws.Range("D108").Value = 3678
$ws.Range("E108").Value = 759

# Georgia (row 148) - updated counters
$ws.Range("B148").Value = 1370
$ws.Range("C148").Value = 9
$ws.Range("D148").Value = 1108
$ws.Range("E148").Value = 245

# Letonia (row 149) - updated counters
$ws.Range("B149").Value = 1327
$ws.Range("C149").Value = 1
$ws.Range("E149").Value = 201

# Birmania (row 173) - updated counters
$ws.Range("B173").Value = 396
$ws.Range("C173").Value = 2
$ws.Range("E173").Value = 57

# Islas Malvinas and Montserrat swap places (rows 213-214)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
